$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.421.93'
$ws.Range('E2').Value = '  -0.34%  '

$ws.Range('D3').Value = '1.924.61'
$ws.Range('E3').Value = '  +3.88%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.10%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.18'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.86%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4739'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '44.36'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.24%  '

$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.2848'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.89%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06591'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.53%  '

$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '19.15'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +8.43%  '

$ws.Range('B12').Value = 'Litecoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '105.11'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +24.36%  '

$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.913.07'
$ws.Range('E13').Value = '  +3.24%  '

$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.07586'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.91%  '

$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.124'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.03%  '

$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.6526'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.33%  '

$ws.Range('B17').Value = 'BitcoinCash'
$ws.Range('C17').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '300.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +21.76%  '

$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '30.420.52'
$ws.Range('E18').Value = '  -0.22%  '

$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.0000'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.04%  '

$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.93%  '

$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.000007516'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.79%  '

$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.166.83'
$ws.Range('E22').Value = '  +2.42%  '

$ws.Range('B23').Value = 'BinanceUSD'
$ws.Range('C23').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9993'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.35%  '

$ws.Range('B24').Value = 'Uniswap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.243'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.24%  '

$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.295'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.39%  '

$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.59'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.40%  '

$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.188'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.81%  '

$ws.Range('B28').Value = 'EthereumClassic'
$ws.Range('C28').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.65'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.32%  '

$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.024'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +8.02%  '

$ws.Range('B30').Value = 'Stellar'
$ws.Range('C30').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.1119'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +9.55%  '

$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.361'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.46%  '

$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.104'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.26%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.920'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.43%  '

$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05005'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +3.45%  '

$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7388'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +5.34%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.141'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.70%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.717'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.56%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01950'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.95%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.696'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.15%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.046'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.53%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8722'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.28%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.24'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.65%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.806'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.69%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9999'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.07%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '69.16'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +10.08%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4124'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.88%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.230'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.185'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.44%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1203'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '34.62'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.27%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05618'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.56%  '
